$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K) values for rows 2-5 as per regenerated save_data
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
